$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.655.86"
$ws.Range("E2").Value = "  -2.20%  "
$ws.Range("D3").Value = "1.761.16"
$ws.Range("E3").Value = "  -3.04%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4309"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3606"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07589"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.112"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.48%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.076"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.236"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "1.760.13"
$ws.Range("E16").Value = "  -3.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001068"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06438"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.54%  "
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.883"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.11%  "
$ws.Range("D23").Value = "27.700.65"
$ws.Range("E23").Value = "  -2.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.085"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("D28").Value = "1.959.68"
$ws.Range("E28").Value = "  -3.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.157"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.100"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.691"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.607"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08951"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02304"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2117"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06019"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6356"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.954"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.74%  "
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.400"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.913"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.63%  "
$ws.Range("E45").Value = "  -4.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5933"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.710"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.991"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.171"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06877"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.80%  "
